$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 147; $r++) {
    $d = $ws.Cells.Item($r, 4)
    $e = $ws.Cells.Item($r, 5)

    $dv = $d.Value2
    $ev = $e.Value2

    if ($dv -like "Ball_*") {
        $d.Value2 = $dv -replace "^Ball_", "Fixation_"
    }
    if ($ev -like "Ball_*") {
        $e.Value2 = $ev -replace "^Ball_", "Fixation_"
    }
}

$ws.Range("F8").Select()
